# Scheduled-runner price/profit refresh across the Leve profit sheets.
# Updates currentAveragePrice*, LevePrice*, and LeveProfit* columns (H-N)
# for the rows whose upstream market data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 91921.55
$ws.Range("I28").Value = 118632.06
$ws.Range("J28").Value = 1105.8
$ws.Range("K28").Value = 118632.06
$ws.Range("L28").Value = 1105.8
$ws.Range("M28").Value = -118147.06
$ws.Range("N28").Value = -2075.8

$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()

$ws.Range("H113").Value = 6680.8
$ws.Range("J113").Value = 6833.3335
$ws.Range("L113").Value = 6833.3335
$ws.Range("N113").Value = -13341.3335

$ws.Range("H129").Value = 2467.15
$ws.Range("J129").Value = 2509.1177
$ws.Range("L129").Value = 7527.353099999999
$ws.Range("N129").Value = -17527.3531

$ws.Range("H138").Value = 3213.182
$ws.Range("I138").Value = 1565.7812
$ws.Range("K138").Value = 4697.3436
$ws.Range("M138").Value = 442.6563999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2470.3333
$ws.Range("I2").Value = 1830.5
$ws.Range("K2").Value = 1830.5
$ws.Range("M2").Value = -1717.5

$ws.Range("H5").Value = 102
$ws.Range("I5").Value = 102
$ws.Range("K5").Value = 102
$ws.Range("M5").Value = 10

$ws.Range("H44").Value = 69946
$ws.Range("J44").Value = 69946
$ws.Range("L44").Value = 69946
$ws.Range("N44").Value = -70922

$ws.Range("H45").Value = 5523.875
$ws.Range("I45").Value = 4598.7144
$ws.Range("J45").Value = 12000
$ws.Range("K45").Value = 4598.7144
$ws.Range("L45").Value = 12000
$ws.Range("M45").Value = -4221.7144
$ws.Range("N45").Value = -12754

$ws.Range("H55").Value = 58287.5
$ws.Range("J55").Value = 63935.4
$ws.Range("L55").Value = 63935.4
$ws.Range("N55").Value = -64565.4

$ws.Range("H116").Value = 2470.3333
$ws.Range("I116").Value = 1830.5
$ws.Range("K116").Value = 1830.5
$ws.Range("M116").Value = 463.5

$ws.Range("H139").Value = 88925.836
$ws.Range("J139").Value = 88925.836
$ws.Range("L139").Value = 88925.836
$ws.Range("N139").Value = -99205.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2470.3333
$ws.Range("I3").Value = 1830.5
$ws.Range("K3").Value = 1830.5
$ws.Range("M3").Value = -1716.5

$ws.Range("H4").Value = 102
$ws.Range("I4").Value = 102
$ws.Range("K4").Value = 102
$ws.Range("M4").Value = 13

$ws.Range("H105").Value = 3271.818
$ws.Range("I105").Value = 2999
$ws.Range("J105").Value = 6000
$ws.Range("K105").Value = 2999
$ws.Range("L105").Value = 6000
$ws.Range("M105").Value = -1252
$ws.Range("N105").Value = -9494

$ws.Range("H107").Value = 1164.5714
$ws.Range("I107").Value = 954.1875
$ws.Range("K107").Value = 954.1875
$ws.Range("M107").Value = 965.8125

$ws.Range("H132").Value = 79999.336
$ws.Range("J132").Value = 79999.336
$ws.Range("L132").Value = 79999.336
$ws.Range("N132").Value = -90119.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 33.636364
$ws.Range("I7").Value = 45.714287
$ws.Range("J7").Value = 12.5
$ws.Range("K7").Value = 45.714287
$ws.Range("L7").Value = 12.5
$ws.Range("M7").Value = 67.285713
$ws.Range("N7").Value = -238.5

$ws.Range("H16").Value = 100000840
$ws.Range("I16").Value = 142857920
$ws.Range("J16").Value = 999.6667
$ws.Range("K16").Value = 142857920
$ws.Range("L16").Value = 999.6667
$ws.Range("M16").Value = -142857633
$ws.Range("N16").Value = -1573.6667

$ws.Range("H31").Value = 1640.2833
$ws.Range("I31").Value = 972
$ws.Range("J31").Value = 2085.8057
$ws.Range("K31").Value = 972
$ws.Range("L31").Value = 2085.8057
$ws.Range("M31").Value = -677
$ws.Range("N31").Value = -2675.8057

$ws.Range("H34").Value = 1640.2833
$ws.Range("I34").Value = 972
$ws.Range("J34").Value = 2085.8057
$ws.Range("K34").Value = 972
$ws.Range("L34").Value = 2085.8057
$ws.Range("M34").Value = -770
$ws.Range("N34").Value = -2489.8057

$ws.Range("H113").Value = 100000840
$ws.Range("I113").Value = 142857920
$ws.Range("J113").Value = 999.6667
$ws.Range("K113").Value = 142857920
$ws.Range("L113").Value = 999.6667
$ws.Range("M113").Value = -142855750
$ws.Range("N113").Value = -5339.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 175264740
$ws.Range("I4").Value = 147845660
$ws.Range("J4").Value = 266661650
$ws.Range("K4").Value = 443536980
$ws.Range("L4").Value = 799984950
$ws.Range("M4").Value = -443536868
$ws.Range("N4").Value = -799985174

$ws.Range("H19").Value = 3786.5
$ws.Range("I19").Value = 2549
$ws.Range("J19").Value = 4199
$ws.Range("K19").Value = 7647
$ws.Range("L19").Value = 12597
$ws.Range("M19").Value = -7473
$ws.Range("N19").Value = -12945

$ws.Range("H107").Value = 1361.8889
$ws.Range("J107").Value = 1625.4117
$ws.Range("L107").Value = 4876.2351
$ws.Range("N107").Value = -8716.2351

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 14814.406
$ws.Range("I113").Value = 1722.4762
$ws.Range("J113").Value = 39808.09
$ws.Range("K113").Value = 1722.4762
$ws.Range("L113").Value = 39808.09
$ws.Range("M113").Value = 447.5237999999999
$ws.Range("N113").Value = -44148.09

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2790.3572
$ws.Range("I61").Value = 1124
$ws.Range("K61").Value = 1124
$ws.Range("M61").Value = -922

$ws.Range("H93").Value = 34483950
$ws.Range("I93").Value = 50001020
$ws.Range("J93").Value = 1577
$ws.Range("K93").Value = 50001020
$ws.Range("L93").Value = 1577
$ws.Range("M93").Value = -49999772
$ws.Range("N93").Value = -4073

$ws.Range("H99").Value = 65910.75
$ws.Range("I99").Value = 52222.832
$ws.Range("K99").Value = 52222.832
$ws.Range("M99").Value = -49227.832

$ws.Range("H113").Value = 2790.3572
$ws.Range("I113").Value = 1124
$ws.Range("K113").Value = 1124
$ws.Range("M113").Value = 1046

$ws.Range("H122").Value = 9138.286
$ws.Range("I122").Value = 5661
$ws.Range("K122").Value = 16983
$ws.Range("M122").Value = -14533

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 31542.5
$ws.Range("I58").Value = 31542.5
$ws.Range("K58").Value = 31542.5
$ws.Range("M58").Value = -31234.5

$ws.Range("H62").Value = 4608.1816
$ws.Range("J62").Value = 5117.5
$ws.Range("L62").Value = 5117.5
$ws.Range("N62").Value = -6365.5

$ws.Range("H65").Value = 4608.1816
$ws.Range("J65").Value = 5117.5
$ws.Range("L65").Value = 25587.5
$ws.Range("N65").Value = -31827.5

$ws.Range("H100").Value = 944.5
$ws.Range("I100").Value = 948.5454999999999
$ws.Range("J100").Value = 900
$ws.Range("K100").Value = 1897.091
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -1356.091
$ws.Range("N100").Value = -2882

$ws.Range("H107").Value = 493.5
$ws.Range("I107").Value = 493.5
$ws.Range("K107").Value = 1480.5
$ws.Range("M107").Value = 439.5

$ws.Range("H136").Value = 47324.652
$ws.Range("I136").Value = 3239.077
$ws.Range("K136").Value = 9717.231
$ws.Range("M136").Value = -7167.231
